$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set column widths for new columns E (5) and F (6), matching A/B widths ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Region / prey species"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "# of stocks"

$ws.Range("E2").Value = "Europe"

$ws.Range("A3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Atlantic herring (Clupea harengus)"
$ws.Range("F3").Value = 12

$ws.Range("A3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "Atlantic herring (Clupea harengus)"
$ws.Range("F4").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "Atlantic mackerel (Scomber scombrus)"
$ws.Range("F5").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "Blue whiting (Micromesistius poutassou)"
$ws.Range("F6").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "Capelin (Mallotus villosus)"
$ws.Range("F7").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "European sprat (Sprattus sprattus)"
$ws.Range("F8").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "Lesser sand-eel (Ammodytes marinus)"
$ws.Range("F9").Value = 3

$ws.Range("A3").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "Lesser sand-eel (Ammodytes spp.)"
$ws.Range("F10").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "Norway pout (Trisopterus esmarkii)"
$ws.Range("F11").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "Whiting (Merlangius merlangus)"
$ws.Range("F12").Value = 3

$ws.Range("E13").Value = "Humboldt Current"

$ws.Range("A3").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = "Peruvian anchoveta (Engraulis ringens)"
$ws.Range("F14").Value = 4

$ws.Range("E15").Value = "South Africa"

$ws.Range("A3").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "European anchovy (Engraulis encrasicolus)"
$ws.Range("F16").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "Pacific sardine (Sardinops sagax)"
$ws.Range("F17").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "Pacific sardine (Sardinops sagax)"
$ws.Range("F18").Value = 1

$ws.Range("E19").Value = "US/Canada East Coast"

$ws.Range("A3").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Atlantic herring (Clupea harengus)"
$ws.Range("F20").Value = 6

$ws.Range("A3").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "Atlantic mackerel (Scomber scombrus)"
$ws.Range("F21").Value = 2

$ws.Range("A3").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "Atlantic menhaden (Brevoortia tyrannus)"
$ws.Range("F22").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = "Gulf menhaden (Brevoortia patronus)"
$ws.Range("F23").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24").Value = "Longfin inshore squid (Doryteuthis pealeii)"
$ws.Range("F24").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = "Northern shortfin squid (Illex illecebrosus)"
$ws.Range("F25").Value = 2

$ws.Range("E26").Value = "US/Canada West Coast"

$ws.Range("A3").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = "Northern anchovy (Engraulis mordax)"
$ws.Range("F27").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "Pacific chub mackerel (Scomber japonicus)"
$ws.Range("F28").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "Pacific hake (Merluccius productus)"
$ws.Range("F29").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "Pacific sardine (Sardinops sagax)"
$ws.Range("F30").Value = 1

$ws.Range("A3").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "Rockfish spp. (Sebastes spp.)"
$ws.Range("F31").Value = 3

$ws.Range("A3").Copy()
$ws.Range("E32").PasteSpecial(-4122)
$ws.Range("E32").Value = "Walleye pollock (Theragra chalcogramma)"
$ws.Range("F32").Value = 4

# --- Clear clipboard / marching ants and set final selection to match target ---
$excel.CutCopyMode = 0
$ws.Range("H11").Select() | Out-Null
